$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing F/G column values (rows 272-435) ---
$ws.Cells.Item(272, 6).Value2 = 30006
$ws.Cells.Item(272, 7).Value2 = 1599
$ws.Cells.Item(273, 6).Value2 = 31572
$ws.Cells.Item(276, 6).Value2 = 11449
$ws.Cells.Item(278, 6).Value2 = 30625
$ws.Cells.Item(285, 6).Value2 = 42361
$ws.Cells.Item(286, 6).Value2 = 55188
$ws.Cells.Item(290, 6).Value2 = 17629
$ws.Cells.Item(291, 6).Value2 = 15166
$ws.Cells.Item(293, 6).Value2 = 83213
$ws.Cells.Item(294, 6).Value2 = 94026
$ws.Cells.Item(295, 6).Value2 = 17301
$ws.Cells.Item(296, 6).Value2 = 1861
$ws.Cells.Item(297, 6).Value2 = 2418
$ws.Cells.Item(298, 6).Value2 = 3267
$ws.Cells.Item(300, 6).Value2 = 72956
$ws.Cells.Item(301, 6).Value2 = 72372
$ws.Cells.Item(303, 6).Value2 = 9680
$ws.Cells.Item(304, 6).Value2 = 6124
$ws.Cells.Item(305, 6).Value2 = 3394
$ws.Cells.Item(306, 6).Value2 = 75684
$ws.Cells.Item(308, 6).Value2 = 15715
$ws.Cells.Item(311, 6).Value2 = 61521
$ws.Cells.Item(312, 6).Value2 = 28209
$ws.Cells.Item(313, 6).Value2 = 76447
$ws.Cells.Item(314, 6).Value2 = 64417
$ws.Cells.Item(315, 6).Value2 = 56426
$ws.Cells.Item(316, 6).Value2 = 50765
$ws.Cells.Item(318, 6).Value2 = 49353
$ws.Cells.Item(319, 6).Value2 = 41397
$ws.Cells.Item(321, 6).Value2 = 90262
$ws.Cells.Item(322, 6).Value2 = 109646
$ws.Cells.Item(323, 6).Value2 = 217598
$ws.Cells.Item(324, 6).Value2 = 250077
$ws.Cells.Item(325, 6).Value2 = 775324
$ws.Cells.Item(325, 7).Value2 = 6504
$ws.Cells.Item(326, 6).Value2 = 418236
$ws.Cells.Item(326, 7).Value2 = 3822
$ws.Cells.Item(327, 6).Value2 = 224078
$ws.Cells.Item(327, 7).Value2 = 2718
$ws.Cells.Item(328, 6).Value2 = 181161
$ws.Cells.Item(333, 6).Value2 = 254988
$ws.Cells.Item(334, 6).Value2 = 193099
$ws.Cells.Item(337, 6).Value2 = 104272
$ws.Cells.Item(337, 7).Value2 = 2923
$ws.Cells.Item(340, 6).Value2 = 387201
$ws.Cells.Item(342, 6).Value2 = 178815
$ws.Cells.Item(343, 6).Value2 = 133569
$ws.Cells.Item(344, 6).Value2 = 135779
$ws.Cells.Item(344, 7).Value2 = 2487
$ws.Cells.Item(345, 6).Value2 = 292276
$ws.Cells.Item(346, 6).Value2 = 674740
$ws.Cells.Item(347, 6).Value2 = 346623
$ws.Cells.Item(349, 6).Value2 = 159756
$ws.Cells.Item(350, 6).Value2 = 127164
$ws.Cells.Item(352, 6).Value2 = 307479
$ws.Cells.Item(354, 6).Value2 = 316764
$ws.Cells.Item(355, 6).Value2 = 222004
$ws.Cells.Item(356, 6).Value2 = 160182
$ws.Cells.Item(357, 6).Value2 = 138219
$ws.Cells.Item(358, 6).Value2 = 158992
$ws.Cells.Item(359, 6).Value2 = 321169
$ws.Cells.Item(363, 6).Value2 = 188516
$ws.Cells.Item(364, 6).Value2 = 168418
$ws.Cells.Item(366, 6).Value2 = 339406
$ws.Cells.Item(369, 6).Value2 = 234733
$ws.Cells.Item(371, 6).Value2 = 160108
$ws.Cells.Item(373, 6).Value2 = 350043
$ws.Cells.Item(375, 6).Value2 = 351315
$ws.Cells.Item(376, 6).Value2 = 221756
$ws.Cells.Item(384, 6).Value2 = 171788
$ws.Cells.Item(389, 6).Value2 = 353694
$ws.Cells.Item(392, 6).Value2 = 220926
$ws.Cells.Item(395, 6).Value2 = 752510
$ws.Cells.Item(401, 6).Value2 = 272679
$ws.Cells.Item(411, 6).Value2 = 225262
$ws.Cells.Item(413, 6).Value2 = 149580
$ws.Cells.Item(422, 6).Value2 = 295000
$ws.Cells.Item(428, 6).Value2 = 100924
$ws.Cells.Item(429, 6).Value2 = 175188
$ws.Cells.Item(429, 7).Value2 = 453
$ws.Cells.Item(430, 6).Value2 = 171549
$ws.Cells.Item(431, 6).Value2 = 164350
$ws.Cells.Item(432, 6).Value2 = 121946
$ws.Cells.Item(432, 7).Value2 = 427
$ws.Cells.Item(433, 6).Value2 = 84994
$ws.Cells.Item(434, 6).Value2 = 78596
$ws.Cells.Item(435, 6).Value2 = 80321
$ws.Cells.Item(435, 7).Value2 = 264

# --- Add new rows 436-438 (data through 17.05.2021) ---
$ws.Cells.Item(436, 1).Value2 = 44330
$ws.Cells.Item(436, 2).Value2 = 387420
$ws.Cells.Item(436, 3).Value2 = 6728
$ws.Cells.Item(436, 4).Value2 = 258
$ws.Cells.Item(436, 5).Value2 = 12203
$ws.Cells.Item(436, 6).Value2 = 135484
$ws.Cells.Item(436, 7).Value2 = 681
$ws.Cells.Item(437, 1).Value2 = 44331
$ws.Cells.Item(437, 2).Value2 = 387523
$ws.Cells.Item(437, 3).Value2 = 2738
$ws.Cells.Item(437, 4).Value2 = 103
$ws.Cells.Item(437, 5).Value2 = 12224
$ws.Cells.Item(437, 6).Value2 = 152761
$ws.Cells.Item(437, 7).Value2 = 1311
$ws.Cells.Item(438, 1).Value2 = 44332
$ws.Cells.Item(438, 2).Value2 = 387659
$ws.Cells.Item(438, 3).Value2 = 3154
$ws.Cells.Item(438, 4).Value2 = 136
$ws.Cells.Item(438, 5).Value2 = 12238
$ws.Cells.Item(438, 6).Value2 = 100939
$ws.Cells.Item(438, 7).Value2 = 195

Write-Output "Updated existing rows and appended rows 436-438."
